$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.016822428083365711
$ws.Range("B1").Value = 0.016822428054804453

$ws.Range("A2").Value = 0.02730670582687315
$ws.Range("B2").Value = -0.027306705851463803

$ws.Range("A3").Value = -0.062737003971333472
$ws.Range("B3").Value = 0.062737003952706608

$ws.Range("A4").Value = -0.012674708612671535
$ws.Range("B4").Value = 0.012674708592601287

$ws.Range("A5").Value = 0.072106449226061797
$ws.Range("B5").Value = -0.072106449248151475
